# Weather workbook refresh: pull latest readings into Sheet 1 and
# re-select the refreshed Temperature column (mirrors the pyexcel re-save).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Activate()

# --- Mumbai ---
$ws1.Cells.Item(2, 2).Value = 26

# --- Delhi (Fahrenheit reading) ---
$ws1.Cells.Item(3, 2).Value = 91.40000000000001

# --- Indore ---
$ws1.Cells.Item(4, 2).Value = 28

# --- Jaipur: no reading available this time, and flag flipped to 0 ---
$ws1.Cells.Item(5, 2).Clear()
$ws1.Cells.Item(5, 4).Value = 0

# --- Bangalore: no reading, unit switched to F, flag flipped to 0 ---
$ws1.Cells.Item(6, 2).Clear()
$ws1.Cells.Item(6, 3).Value = "F"
$ws1.Cells.Item(6, 4).Value = 0

# --- Chennai ---
$ws1.Cells.Item(7, 2).Value = 32

# --- Hyderabad ---
$ws1.Cells.Item(8, 2).Value = 27.74000000000001

# --- Mysore ---
$ws1.Cells.Item(9, 2).Value = 26

# --- Gangtok: unit switched to F ---
$ws1.Cells.Item(10, 2).Value = 76.40600000000003
$ws1.Cells.Item(10, 3).Value = "F"

# --- Guwhati: reading still unavailable ---
$ws1.Cells.Item(11, 2).Value = "NA"

# --- Dubai: no reading available this time ---
$ws1.Cells.Item(12, 2).Clear()

# Shrink the row heights to the tighter auto value used after the refresh
# (header row keeps its original height).
for ($r = 2; $r -le 12; $r++) {
  $ws1.Rows.Item($r).RowHeight = 13.8
}

# Reflect the refreshed column as the active selection (Sheet 2's own
# selection stays at A1, then Sheet 1 is left as the active/visible tab).
$ws2.Range("A1").Select() | Out-Null
$ws1.Range("B2:B12").Select() | Out-Null
